# Loan RBI, Variable Instalments
#
# The "Repayment Schedule" sheet gains a new (blank) column between the
# existing "In Advance" (M) and "Late" (old N) columns, so "Late" moves to
# O and "Outstanding" (old P) moves to Q. "Repayment Schedule" also becomes
# the active/selected sheet (instead of "NewLoanInput"), with a new
# selection on it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column at N - this shifts the old N ("Late") column to
# O and the old P ("Outstanding") column to Q, exactly matching the diff.
$ws.Columns("N").Insert() | Out-Null

# Give the freshly inserted column a sensible (manually-set, non bestFit)
# width, close to its neighbouring "In Advance"/"Late" columns.
$ws.Columns("N").ColumnWidth = 10.3

# Make "Repayment Schedule" the active sheet with the new selection, moving
# the previously active sheet ("NewLoanInput") out of focus.
$ws.Activate() | Out-Null
$ws.Range("U13").Select() | Out-Null
